# formula: add more index functions
#   - HLOOKUP
#   - LOOKUP
#   - TRANSPOSE
#   - VLOOKUP
#
# Applies the edit to the "Indexing" worksheet (adds VLOOKUP / LOOKUP /
# HLOOKUP / TRANSPOSE example formulas and their backing data), and makes
# "Indexing" the active sheet/tab instead of "Statistics".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # "Indexing"
$stats = $wb.Worksheets.Item(6) # "Statistics"

# --- New header labels on row 4 (bold + centered, like F4/G4/H4) -----------
$ws.Range("I4").Value = "VLOOKUP"
$ws.Range("J4").Value = "LOOKUP"
$ws.Range("K4").Value = "HLOOKUP"
$ws.Range("I4:K4").Font.Bold = $true
$ws.Range("I4:K4").HorizontalAlignment = -4108   # xlCenter

# --- VLOOKUP / LOOKUP / HLOOKUP example formulas ---------------------------
$ws.Range("I5").Formula = "=VLOOKUP()"
$ws.Range("J5").Formula = "=LOOKUP()"
$ws.Range("K5").Formula = "=HLOOKUP()"

$ws.Range("I6").Formula = "=VLOOKUP(1,C18:E24)"
$ws.Range("J6").Formula = "=LOOKUP(1)"
$ws.Range("K6").Formula = "=HLOOKUP(1)"

$ws.Range("I7").Formula = "=VLOOKUP(1,$C$18:$E$24,1)"
$ws.Range("J7").Formula = "=LOOKUP(1,C18:C24)"
$ws.Range("K7").Formula = "=HLOOKUP(1,$C$27:$F$29,1)"

$ws.Range("I8").Formula = "=VLOOKUP(1,$C$18:$E$24,2)"
$ws.Range("J8").Formula = "=LOOKUP(2,C18:E18)"
$ws.Range("K8").Formula = "=HLOOKUP(1,$C$27:$F$29,2)"

$ws.Range("I9").Formula = "=VLOOKUP(1,$C$18:$E$24,3)"
$ws.Range("J9").Formula = "=LOOKUP(2,C19:E19)"
$ws.Range("K9").Formula = "=HLOOKUP(1,$C$27:$F$29,3)"

$ws.Range("I10").Formula = "=VLOOKUP(2.3,$C$18:$E$24,1,0)"
$ws.Range("J10").Formula = "=LOOKUP(1,C18:C24,D18:D24)"
$ws.Range("K10").Formula = "=HLOOKUP(1.5,$C$27:$F$29,1)"

$ws.Range("I11").Formula = "=VLOOKUP(2.3,$C$18:$E$24,1,TRUE())"
$ws.Range("J11").Formula = "=LOOKUP(3,C18:C24,E18:E24)"
$ws.Range("K11").Formula = "=HLOOKUP(1.5,$C$27:$F$29,1,1)"

$ws.Range("I12").Formula = "=VLOOKUP(4.2,C18:C24,1)"
$ws.Range("J12").Formula = "=LOOKUP(""A"",C18:E18,C19:E19)"
$ws.Range("K12").Formula = "=HLOOKUP(1.5,$C$27:$F$29,1,0)"

$ws.Range("I13").Formula = "=VLOOKUP(4.2,C19:C25,1,FALSE())"

$ws.Range("I15").Formula = "=VLOOKUP(10,$C$18:$E$24,1)"

$ws.Range("B16").Value = "B13"
$ws.Range("I16").Formula = "=VLOOKUP(-1,$C$18:$E$24,1)"

$ws.Range("I17").Formula = "=VLOOKUP(10,$C$18:$E$24,2)"

$ws.Range("I18").Formula = "=VLOOKUP(10,$C$18:$E$24,3)"
$ws.Range("I19").Formula = "=VLOOKUP(10,$C$18:$E$24,4)"
$ws.Range("I20").Formula = "=VLOOKUP(10,$C$18:$E$24,1)"
$ws.Range("I21").Formula = "=VLOOKUP(10,$C$18:$E$24,0)"
$ws.Range("I22").Formula = "=VLOOKUP(""CC"",D18:E24,1)"
$ws.Range("I23").Formula = "=VLOOKUP(""CC"",D19:E25,2)"

# --- Backing data table used by VLOOKUP (C18:E24) --------------------------
$colC = @(1,2,3,4,5,6,7)
$colD = @("a","b","c","d","e","f","g")
$colE = @("h","I","j","k","l","m","n")
for ($i = 0; $i -lt 7; $i++) {
    $r = 18 + $i
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
    $ws.Cells.Item($r, 5).Value = $colE[$i]
}

# --- Backing data table used by HLOOKUP (C27:F29) --------------------------
$row27 = @(1,2,3,4)
$row28 = @("a","b","c","d")
$row29 = @("e","f","g","h")
for ($i = 0; $i -lt 4; $i++) {
    $c = 3 + $i
    $ws.Cells.Item(27, $c).Value = $row27[$i]
    $ws.Cells.Item(28, $c).Value = $row28[$i]
    $ws.Cells.Item(29, $c).Value = $row29[$i]
}

# --- TRANSPOSE section -------------------------------------------------------
$ws.Range("C32").Value = "TRANSPOSE"
$ws.Range("C32").Font.Bold = $true

$ws.Range("C33:I33").FormulaArray = "=TRANSPOSE(C18:C24)"
$ws.Range("C36:E38").FormulaArray = "=TRANSPOSE(C27:E29)"

$ws.Range("C40").Formula = "=C38"
$ws.Range("C41").Formula = "=D37"
$ws.Range("C42").Formula = "=D33"

# --- Make "Indexing" the active sheet/tab instead of "Statistics" ----------
$ws.Activate()
$ws.Range("C43").Select()
